# Update countries & provincias Spain
# Applies the COVID dashboard data refresh described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "last refreshed" timestamp banner (A1) ---------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 10:35"

# --- 2. Polonia (row 40) gets refreshed totals ---------------------------
$ws.Range("B40").Value = 23987
$ws.Range("C40").Value = 201
$ws.Range("E40").Value = 11473
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 1065

# --- 3. Bulgaria / El Salvador swap rank (rows 86-87) ---------------------
# El Salvador's updated numbers now outrank Bulgaria's, so the two
# countries trade places in the row-86 / row-87 slots while keeping the
# sheet sorted by total cases (column B) descending.
$ws.Range("A86").Value = "El Salvador"
$ws.Range("B86").Value = 2582
$ws.Range("C86").Value = 65
$ws.Range("D86").Value = 1063
$ws.Range("E86").Value = 1473
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 46

$ws.Range("A87").Value = "Bulgaria"
$ws.Range("B87").Value = 2519
$ws.Range("C87").Value = 6
$ws.Range("D87").Value = 1090
$ws.Range("E87").Value = 1289
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 140
